$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay text (avoids Excel auto-numeric coercion)
# without leaving any NumberFormat/style residue on the destination cell.
# Uses a scratch cell (Text-formatted) + copy/PasteSpecial(values-only), then
# fully clears the scratch cell so no stray cell/format is left behind.
$scratch = $ws.Range("ZZ1")

function Set-TextValue([string]$addr, [string]$val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

# --- Direct text/string assignments (values that Excel will not mis-type) ---
$ws.Range("D2").Value = '24.508.06'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '1.659.54'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  -3.27%  '
$ws.Range("E8").Value = '  -3.31%  '
$ws.Range("E9").Value = '  -5.42%  '
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E11").Value = '  -6.05%  '
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("E13").Value = '  -5.44%  '
$ws.Range("E14").Value = '  -6.95%  '
$ws.Range("D15").Value = '1.656.41'
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("E16").Value = '  -5.47%  '
$ws.Range("E17").Value = '  -6.87%  '
$ws.Range("E18").Value = '  -2.52%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  -8.94%  '
$ws.Range("E21").Value = '  -6.26%  '
$ws.Range("E22").Value = '  -8.40%  '
$ws.Range("E23").Value = '  -3.20%  '
$ws.Range("D24").Value = '24.473.68'
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  -16.02%  '
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("E28").Value = '  -8.28%  '
$ws.Range("D29").Value = '1.838.73'
$ws.Range("E29").Value = '  -2.92%  '
$ws.Range("E30").Value = '  -5.21%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("E32").Value = '  -4.57%  '
$ws.Range("E33").Value = '  -16.39%  '
$ws.Range("E34").Value = '  -5.52%  '
$ws.Range("E35").Value = '  -4.90%  '
$ws.Range("E36").Value = '  -9.31%  '
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("E38").Value = '  -7.61%  '
$ws.Range("E39").Value = '  -7.40%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E40").Value = '  -7.63%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E41").Value = '  -8.18%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E42").Value = '  -5.65%  '
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  -8.06%  '
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("E46").Value = '  -8.95%  '
$ws.Range("E47").Value = '  -8.32%  '
$ws.Range("E48").Value = '  -5.74%  '
$ws.Range("E49").Value = '  -8.18%  '
$ws.Range("E50").Value = '  -5.20%  '
$ws.Range("E51").Value = '  -6.09%  '

# --- Values that look numeric and must be forced to stay text ---
Set-TextValue "D4" '1.004'
Set-TextValue "D5" '307.61'
Set-TextValue "D6" '1.001'
Set-TextValue "D7" '0.3619'
Set-TextValue "D8" '47.54'
Set-TextValue "D9" '0.3253'
Set-TextValue "D10" '1.122'
Set-TextValue "D11" '0.07007'
Set-TextValue "D13" '5.896'
Set-TextValue "D14" '19.44'
Set-TextValue "D16" '6.571'
Set-TextValue "D17" '0.00001047'
Set-TextValue "D18" '0.06542'
Set-TextValue "D19" '1.001'
Set-TextValue "D20" '76.48'
Set-TextValue "D21" '5.923'
Set-TextValue "D22" '15.67'
Set-TextValue "D25" '2.466'
Set-TextValue "D26" '2.322'
Set-TextValue "D27" '147.12'
Set-TextValue "D28" '18.47'
Set-TextValue "D30" '124.04'
Set-TextValue "D31" '1.168'
Set-TextValue "D32" '3.980'
Set-TextValue "D33" '5.647'
Set-TextValue "D34" '1.695'
Set-TextValue "D35" '0.08389'
Set-TextValue "D36" '12.36'
Set-TextValue "D37" '5.194'
Set-TextValue "D38" '0.06063'
Set-TextValue "D39" '0.02203'
Set-TextValue "D40" '0.2056'
Set-TextValue "D41" '8.223'
Set-TextValue "D42" '1.203'
Set-TextValue "D44" '0.5898'
Set-TextValue "D45" '3.739'
Set-TextValue "D46" '12.58'
Set-TextValue "D47" '0.5593'
Set-TextValue "D48" '122.18'
Set-TextValue "D49" '1.936'
Set-TextValue "D50" '0.06912'
Set-TextValue "D51" '74.28'

# Clean up the scratch cell completely so it leaves no trace in the sheet
$scratch.ClearContents()
$scratch.ClearFormats()
$scratch.Clear()
$ws.Range("A1").Select() | Out-Null

Write-Host "Applied cryptos list update"
